# Apply the edit described by the diff:
#  1. Remove the (empty) INNING_NUMBER values that were sitting in B9 / B10
#     of the "ODI Batting" sheet.
#  2. Add a new worksheet "ODI Batting Extra" (sheetId 4) after "ODI Bowling"
#     and populate it with the MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#     PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH table.

$wb = $excel.ActiveWorkbook

# --- 1. Clear the stray empty cells in "ODI Batting" -----------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("B9").ClearContents()
$battingSheet.Range("B10").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet at the end -------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bowlingSheet)
$newSheet.Name = "ODI Batting Extra"

# Copy the header formatting (bold, centered, bordered) from an existing
# sheet's header row so the new header cells reuse the same style.
$battingSheet.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

# Header row
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$rows = @(
    @("4086", 7,    "0", "0", "0.54%",  "NO"),
    @("4182", 7,    "3", "0", "7.14%",  "NO"),
    @("4183", $null, $null, $null, $null, "NO"),
    @("4186", 8,    "0", "0", "0.55%",  "NO"),
    @("4187", 8,    "1", "0", "2.14%",  "NO"),
    @("4188", $null, $null, $null, $null, "NO"),
    @("4206", 7,    "2", "0", "11.76%", "NO"),
    @("4247", $null, $null, $null, $null, "NO"),
    @("4261", 6,    $null, $null, $null, "NO"),
    @("4264", 5,    "4", "0", "6.77%",  "NO"),
    @("4488", 7,    "0", "0", "1.77%",  "NO"),
    @("4491", 6,    "0", "0", "1.60%",  "NO")
)

$r = 2
foreach ($row in $rows) {
    # MATCH_CODE - stored as text (digits only)
    $cellA = $newSheet.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]

    # BATTING_POSITION - numeric, or blank
    $cellB = $newSheet.Cells.Item($r, 2)
    if ($null -ne $row[1]) {
        $cellB.Value = $row[1]
    } else {
        $cellB.Value = ""
    }

    # NUM_4 - text, or blank
    $cellC = $newSheet.Cells.Item($r, 3)
    if ($null -ne $row[2]) {
        $cellC.NumberFormat = "@"
        $cellC.Value = $row[2]
    } else {
        $cellC.Value = ""
    }

    # NUM_6 - text, or blank
    $cellD = $newSheet.Cells.Item($r, 4)
    if ($null -ne $row[3]) {
        $cellD.NumberFormat = "@"
        $cellD.Value = $row[3]
    } else {
        $cellD.Value = ""
    }

    # PERCENT_RUNS_OF_TOTAL - text, or blank
    $cellE = $newSheet.Cells.Item($r, 5)
    if ($null -ne $row[4]) {
        $cellE.NumberFormat = "@"
        $cellE.Value = $row[4]
    } else {
        $cellE.Value = ""
    }

    # MAN_OF_MATCH - text
    $newSheet.Cells.Item($r, 6).Value = $row[5]

    $r++
}

$newSheet.Range("A1").Select() | Out-Null
